# regen sval data to filter save games
# Update the numeric stat columns (B:E, G) for rows 2-7 on the active sheet
# with the recomputed values from the filtered save-game dataset.
# Column F ("Win") values are unchanged by this regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 4.371470058157054

$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 1.896700893398075

$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 3.055818435266709

$ws.Range("B5").Value = 3.182878228561681
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 5.488907176552729

$ws.Range("B6").Value = 1.505614041169197
$ws.Range("C6").Value = 9.226618575922256
$ws.Range("D6").Value = 3.082599426703578
$ws.Range("E6").Value = 6.48142807727062
$ws.Range("G6").Value = 20.29626012106565

$ws.Range("B7").Value = 3.182878228561681
$ws.Range("C7").Value = 1.65323645889881
$ws.Range("D7").Value = 0.7127328510149897
$ws.Range("E7").Value = 0.4998867070740569
$ws.Range("G7").Value = 6.048734245549538
